# Update the "Date" column (B2:B17) on the NumberError sheet with the
# new execution timestamps recorded for the latest RAD BeforePayments
# test run (Mon Sep 11 2023 run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value  = "Mon Sep 11 13:53:45 EDT 2023"
$ws.Range("B3").Value  = "Mon Sep 11 13:53:54 EDT 2023"
$ws.Range("B4").Value  = "Mon Sep 11 13:54:04 EDT 2023"
$ws.Range("B5").Value  = "Mon Sep 11 13:54:14 EDT 2023"
$ws.Range("B6").Value  = "Mon Sep 11 13:54:24 EDT 2023"
$ws.Range("B7").Value  = "Mon Sep 11 13:54:33 EDT 2023"
$ws.Range("B8").Value  = "Mon Sep 11 13:54:43 EDT 2023"
$ws.Range("B9").Value  = "Mon Sep 11 13:54:53 EDT 2023"
$ws.Range("B10").Value = "Mon Sep 11 13:55:03 EDT 2023"
$ws.Range("B11").Value = "Mon Sep 11 13:55:13 EDT 2023"
$ws.Range("B12").Value = "Mon Sep 11 13:55:22 EDT 2023"
$ws.Range("B13").Value = "Mon Sep 11 13:55:32 EDT 2023"
$ws.Range("B14").Value = "Mon Sep 11 13:55:42 EDT 2023"
$ws.Range("B15").Value = "Mon Sep 11 13:55:52 EDT 2023"
$ws.Range("B16").Value = "Mon Sep 11 13:56:02 EDT 2023"
$ws.Range("B17").Value = "Mon Sep 11 13:56:12 EDT 2023"
